# Add a new row (row 89) of match data to the worksheet, mirroring the
# structure of the existing rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 89

$ws.Cells.Item($row, 1).Value = 88
$ws.Cells.Item($row, 2).Value = "montenegro"
$ws.Cells.Item($row, 3).Value = "prva-crnogorska-liga"
$ws.Cells.Item($row, 4).Value = "2023-2024"
$ws.Cells.Item($row, 5).Value = 45262.64583333334
$ws.Cells.Item($row, 6).Value = "Petrovac"
$ws.Cells.Item($row, 7).Value = 2
$ws.Cells.Item($row, 8).Value = "Buducnost"
$ws.Cells.Item($row, 9).Value = 2
$ws.Cells.Item($row, 10).Value = 4.11
$ws.Cells.Item($row, 11).Value = "01/12/2023 03:43"
$ws.Cells.Item($row, 12).Value = 5.31
$ws.Cells.Item($row, 13).Value = "02/12/2023 15:13"
$ws.Cells.Item($row, 14).Value = 3.37
$ws.Cells.Item($row, 15).Value = "01/12/2023 03:43"
$ws.Cells.Item($row, 16).Value = 3.85
$ws.Cells.Item($row, 17).Value = "02/12/2023 15:13"
$ws.Cells.Item($row, 18).Value = 1.76
$ws.Cells.Item($row, 19).Value = "01/12/2023 03:43"
$ws.Cells.Item($row, 20).Value = 1.61
$ws.Cells.Item($row, 21).Value = "02/12/2023 15:13"
$ws.Cells.Item($row, 22).Value = "https://www.betexplorer.com/football/montenegro/prva-crnogorska-liga/petrovac-buducnost/M7xM4XSj/"

# Copy styles from row 88 to keep formatting (bold/border on col A, date format on col E)
$ws.Range("A88").Copy()
$ws.Range("A$row").PasteSpecial(-4122) # xlPasteFormats

$ws.Range("E88").Copy()
$ws.Range("E$row").PasteSpecial(-4122) # xlPasteFormats
